$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '34.620.73'
$ws.Range("E2").Value = '  +2.08%  '

# Row 3
$ws.Range("D3").Value = '1.854.67'
$ws.Range("E3").Value = '  +4.75%  '

# Row 5
$ws.Range("D5").Value = '226.93'
$ws.Range("E5").Value = '  +0.83%  '

# Row 6
$ws.Range("D6").Value = '0.557'
$ws.Range("E6").Value = '  +2.49%  '

# Row 7
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.58%  '

# Row 8
$ws.Range("D8").Value = '32.79'
$ws.Range("E8").Value = '  +5.52%  '

# Row 9
$ws.Range("D9").Value = '0.296'
$ws.Range("E9").Value = '  +6.32%  '

# Row 10
$ws.Range("D10").Value = '0.0720'
$ws.Range("E10").Value = '  +10.24%  '

# Row 11
$ws.Range("E11").Value = '  +0.78%  '

# Row 12
$ws.Range("D12").Value = '2.110.41'
$ws.Range("E12").Value = '  +4.10%  '

# Row 13
$ws.Range("D13").Value = '1.853.08'
$ws.Range("E13").Value = '  +4.77%  '

# Row 14
$ws.Range("D14").Value = '11.26'
$ws.Range("E14").Value = '  +3.64%  '

# Row 15
$ws.Range("D15").Value = '0.656'
$ws.Range("E15").Value = '  +6.05%  '

# Row 16
$ws.Range("D16").Value = '34.639.22'
$ws.Range("E16").Value = '  +2.09%  '

# Row 17
$ws.Range("D17").Value = '4.39'
$ws.Range("E17").Value = '  +5.16%  '

# Row 18
$ws.Range("D18").Value = '70.13'
$ws.Range("E18").Value = '  +2.87%  '

# Row 19
$ws.Range("D19").Value = '254.62'
$ws.Range("E19").Value = '  +1.77%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0811'
$ws.Range("E20").Value = '  +10.66%  '

# Row 21
$ws.Range("E21").Value = '  +10.88%  '

# Row 22
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.46%  '

# Row 23
$ws.Range("D23").Value = '4.35'
$ws.Range("E23").Value = '  +4.77%  '

# Row 24
$ws.Range("E24").Value = '  +1.40%  '

# Row 25
$ws.Range("D25").Value = '161.83'
$ws.Range("E25").Value = '  +4.23%  '

# Row 26
$ws.Range("D26").Value = '17.03'
$ws.Range("E26").Value = '  +4.24%  '

# Row 27
$ws.Range("D27").Value = '7.30'
$ws.Range("E27").Value = '  +5.25%  '

# Row 28
$ws.Range("E28").Value = '  +2.53%  '

# Row 29
$ws.Range("D29").Value = '0.997'

# Row 30
$ws.Range("D30").Value = '0.0537'
$ws.Range("E30").Value = '  +5.44%  '

# Row 31
$ws.Range("E31").Value = '  +2.92%  '

# Row 32
$ws.Range("D32").Value = '1.22'
$ws.Range("E32").Value = '  +2.35%  '

# Row 33
$ws.Range("D33").Value = '515.83'
$ws.Range("E33").Value = '  +890.98%  '

# Row 34
$ws.Range("D34").Value = '3.66'
$ws.Range("E34").Value = '  +3.63%  '

# Row 35
$ws.Range("D35").Value = '1.97'
$ws.Range("E35").Value = '  +8.25%  '

# Row 36
$ws.Range("D36").Value = '1.457.57'
$ws.Range("E36").Value = '  +1.11%  '

# Row 37
$ws.Range("D37").Value = '0.662'
$ws.Range("E37").Value = '  +6.58%  '

# Row 38
$ws.Range("D38").Value = '1.08'
$ws.Range("E38").Value = '  +2.51%  '

# Row 39
$ws.Range("D39").Value = '0.0195'
$ws.Range("E39").Value = '  +5.51%  '

# Row 40
$ws.Range("D40").Value = '0.989'
$ws.Range("E40").Value = '  +12.42%  '

# Row 41
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.88'
$ws.Range("E41").Value = '  +1.20%  '

# Row 42
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '83.69'
$ws.Range("E42").Value = '  +1.94%  '

# Row 43
$ws.Range("E43").Value = '  +0.87%  '

# Row 44
$ws.Range("E44").Value = '  +6.87%  '

# Row 45
$ws.Range("E45").Value = '  +6.99%  '

# Row 46
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '2.008.12'
$ws.Range("E46").Value = '  +4.25%  '

# Row 47
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = '12.61'
$ws.Range("E47").Value = '  +7.79%  '

# Row 48
$ws.Range("E48").Value = '  +0.09%  '

# Row 49
$ws.Range("D49").Value = '0.0494'
$ws.Range("E49").Value = '  -2.71%  '

# Row 50
$ws.Range("D50").Value = '106.97'
$ws.Range("E50").Value = '  +10.47%  '

# Row 51
$ws.Range("E51").Value = '  -0.04%  '
